$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# --- Move row 7 ("5. Use the text file...") down to row 10, leaving old row 7 empty ---
$row7Text = $ws.Range("A7").Text
$ws.Range("A10").Value = $row7Text
$ws.Range("A7").Clear()

# --- Rows 7-9: the three save-format bullet points (col B), quote-prefixed since they start with '-' ---
$ws.Range("B7").Value = "'- ""Text (tab-delimited) (*.txt)"""
$ws.Range("B8").Value = "'- ""CSV (comma-delimited) (*.csv)"""
$ws.Range("B9").Value = "'- Excel files (*.xlsx and *.xls)"

# --- Row 6: replace the old combined "4. ..." instruction with the new shorter lead-in ---
$ws.Range("A6").Value = "4. Save as the excel file (only the ""Fill out this form"" sheet) as any of the followings:"

# --- Column A is much narrower now that the long instruction text wraps across rows ---
$ws.Columns.Item(1).ColumnWidth = 11.28515625

# --- Selection moves back up to A2 ---
[void]$ws.Range("A2").Select()
